# adding new progress as of date 04 nov 2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Scratch cell (outside the used range) formatted as Text, used as a
# staging area so typed-looking date strings ("04-Nov-2025", ...) are
# written as literal text instead of being auto-converted into real
# date serial numbers.
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"

function Set-TextValue($cellAddr, $text, $formatSourceAddr) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)   # xlPasteValues
    $ws.Range($formatSourceAddr).Copy()
    $ws.Range($cellAddr).PasteSpecial(-4122)   # xlPasteFormats
}

# --- Row 3: period-to-expire / last-update refresh ---
$ws.Range("H3").Value = 534
Set-TextValue "I3" "04-Nov-2025" "I3"

# --- Row 4: period-to-expire / last-update refresh ---
$ws.Range("H4").Value = 534
Set-TextValue "I4" "04-Nov-2025" "I4"

# --- Row 5: training got completed, record becomes VALID ---
# Copy the normal "valid" row formatting (row 3) onto row 5 so the
# special NOT-VALID highlight style is no longer used.
$ws.Range("A3:K3").Copy()
$ws.Range("A5:K5").PasteSpecial(-4122)   # xlPasteFormats

Set-TextValue "F5" "23-Apr-2025" "F3"
Set-TextValue "G5" "23-Apr-2027" "G3"
$ws.Range("H5").Value = 534
Set-TextValue "I5" "04-Nov-2025" "I3"
$ws.Range("J5").Value = "VALID"
$ws.Range("K5").Value = ""

# --- Row 6: period-to-expire / last-update refresh ---
$ws.Range("H6").Value = 604
Set-TextValue "I6" "04-Nov-2025" "I6"

# Remove the scratch cell so it doesn't show up as extra used range.
$scratch.Clear()

# --- Column width tweaks (J & K got narrower) ---
# (ColumnWidth uses "characters" units that round-trip to the stored
# OOXML <col width> through Excel's nonlinear pixel-snapping formula;
# 7.15 / 8.15 are the values that land exactly on stored widths 8 / 9.)
$ws.Range("J1").EntireColumn.ColumnWidth = 7.15
$ws.Range("K1").EntireColumn.ColumnWidth = 8.15
